$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark row 68 (last existing detection point row) as the 60th ground truth point
$ws.Cells.Item(68, 5).Value = "60th Ground Truth Point"

# Append new accession numbers (ground truth points without detection data yet)
$ws.Cells.Item(69, 1).Value = 8609024
$ws.Cells.Item(70, 1).Value = 8649140
$ws.Cells.Item(71, 1).Value = 8664266
$ws.Cells.Item(72, 1).Value = 8666303
$ws.Cells.Item(73, 1).Value = 8687072
$ws.Cells.Item(74, 1).Value = 8697536
$ws.Cells.Item(75, 1).Value = 8717421
$ws.Cells.Item(76, 1).Value = 8723551
$ws.Cells.Item(77, 1).Value = 8734298
$ws.Cells.Item(78, 1).Value = 8766117
$ws.Cells.Item(78, 5).Value = "70th Ground Truth Point"
$ws.Cells.Item(79, 1).Value = 8791550
$ws.Cells.Item(80, 1).Value = 8837927
$ws.Cells.Item(81, 1).Value = 8847103
$ws.Cells.Item(82, 1).Value = 8864910
$ws.Cells.Item(83, 1).Value = 8915108
$ws.Cells.Item(84, 1).Value = 8920727
$ws.Cells.Item(85, 1).Value = 8931170
$ws.Cells.Item(86, 1).Value = 8931305
$ws.Cells.Item(87, 1).Value = 8943923
$ws.Cells.Item(88, 1).Value = 8958140
$ws.Cells.Item(88, 5).Value = "80th Ground Truth Point"
$ws.Cells.Item(89, 1).Value = 8976584
$ws.Cells.Item(90, 1).Value = 8985647
$ws.Cells.Item(91, 1).Value = 9014200
$ws.Cells.Item(92, 1).Value = 9020776
$ws.Cells.Item(93, 1).Value = 9043600
$ws.Cells.Item(94, 1).Value = 9049401
$ws.Cells.Item(95, 1).Value = 9063971
$ws.Cells.Item(96, 1).Value = 9086188
$ws.Cells.Item(97, 1).Value = 9104767
$ws.Cells.Item(98, 1).Value = 9115352
$ws.Cells.Item(98, 5).Value = "90th Ground Truth Point"
$ws.Cells.Item(99, 1).Value = 9118382
$ws.Cells.Item(100, 1).Value = 9136320
$ws.Cells.Item(101, 1).Value = 9145963
$ws.Cells.Item(102, 1).Value = 9189512

# Update the active selection to reflect the new end of the appended data
$ws.Range("B69").Select() | Out-Null
